$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "Tortilla" row (row 10) entirely
$ws.Rows("10").Delete()

# Fix the typo in the Pistachios item code: 836067 -> 835067
$ws.Range("A8").Value = "835067 PISTACHIOS SALTED"
